$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ConstValue")

# 1) Add the new config row: ChatGetCount / int / 10 / 聊天消息单次拉取数量
#    Copy row 8's formatting down to row 9 first, then fill in the new values.
$ws.Range("C8:F8").Copy($ws.Range("C9:F9")) | Out-Null
$ws.Range("C9").Value = "ChatGetCount"
$ws.Range("D9").Value = "int"
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = "聊天消息单次拉取数量"

# 2) Correct D8's declared type from "int" to "string" -- copy D6 (which already
#    holds the "string" value with the matching cell style) down onto D8 so both
#    the value and formatting move together.
$ws.Range("D6").Copy($ws.Range("D8")) | Out-Null

# 3) Rows no longer need an explicit height override once re-laid-out.
$ws.Rows("2:9").AutoFit() | Out-Null

# 4) Move the active selection the way it ends up after typing into F9 and
#    pressing Enter / moving down a row.
$ws.Range("C10").Select() | Out-Null

$wb.Save()
